# Update the build timestamp embedded in the version string throughout the
# workbook: "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Komsomolskaya Coal Mine, Russia, M2339, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S ("build_version") holds the same version string on every data
# row (rows 2 through the last used row); update them all.
$lastRow = $data.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $data.Cells.Item($r, 19).Value = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
}
